$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pairs of rows whose data (columns B:AC) got swapped between them,
# while column A (the running "id" index) stays fixed on its own row.
$swapPairs = @(
    @(17, 18),
    @(54, 55),
    @(73, 74),
    @(103, 104),
    @(135, 136),
    @(161, 162)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")
    $vals1 = $range1.Value2
    $vals2 = $range2.Value2
    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}

# Append the new fixture as row 165, matching the formatting used for
# the "id" column (A) and the "Date" column (E) on the previous row.
$ws.Range("A164").Copy()
$ws.Range("A165").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E164").Copy()
$ws.Range("E165").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A165").Value = 163
$ws.Range("B165").Value = 8122430
$ws.Range("C165").Value = "Iraq League"
$ws.Range("D165").Value = "Iraq League"
$ws.Range("E165").Value2 = 45403.51041666666
$ws.Range("F165").Value = "Naft Maysan"
$ws.Range("G165").Value = "Al Kahrabaa"
$ws.Range("K165").Value = 2.15
$ws.Range("L165").Value = 2.75
$ws.Range("M165").Value = 3.5
$ws.Range("N165").Value = 2.15
$ws.Range("O165").Value = 2.75
$ws.Range("P165").Value = 3.5
$ws.Range("Q165").Value = -0.25
$ws.Range("R165").Value = 1.875
$ws.Range("S165").Value = 1.925
$ws.Range("T165").Value = 2
$ws.Range("U165").Value = 2
$ws.Range("V165").Value = 1.8
$ws.Range("W165").Value = 0
$ws.Range("X165").Value = 0
$ws.Range("Y165").Value = 0
$ws.Range("Z165").Value = 0
$ws.Range("AA165").Value = 0
